$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5005.48337366473
$ws.Range("B3").Value = 4922.298456099087
$ws.Range("B4").Value = 4850.625419333438
$ws.Range("B5").Value = 4797.274898040722
$ws.Range("B6").Value = 4859.095012476347
$ws.Range("B7").Value = 4960.889853826959
$ws.Range("B8").Value = 5211.449593293681
$ws.Range("B9").Value = 6049.649550323782
$ws.Range("B10").Value = 8417.127547097272
$ws.Range("B11").Value = 13902.175
$ws.Range("B12").Value = 15210.828
$ws.Range("B13").Value = 15097.205
$ws.Range("B14").Value = 14468.754
$ws.Range("B15").Value = 14615.9345
$ws.Range("B16").Value = 15642.9465
$ws.Range("B17").Value = 15864.0115
$ws.Range("B18").Value = 16275.092
$ws.Range("B19").Value = 16265.2175
$ws.Range("B20").Value = 15390.046
$ws.Range("B21").Value = 13874.261
$ws.Range("B22").Value = 11951.5345
$ws.Range("B23").Value = 9041.0965
$ws.Range("B24").Value = 6150.5125
$ws.Range("B25").Value = 5731.305
